$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "https://boards.greenhouse.io/lambda/jobs/6053055003"
$ws.Range("B2").Value = "Remote but location not found"

# Delete row 3 entirely (shifts rows up, removing the row)
$ws.Rows.Item(3).Delete()
